$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.039.65'
$ws.Range("E2").Value = '  -0.05%  '

$ws.Range("D3").Value = '1.833.32'
$ws.Range("E3").Value = '  -0.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9954'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.40%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6242'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9981'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07488'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.85%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2939'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.32'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07694'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.72%  '

$ws.Range("D12").Value = '1.845.33'
$ws.Range("E12").Value = '  +0.56%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.023'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6747'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.75'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.30%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009378'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.57%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.970'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.61%  '

$ws.Range("D18").Value = '29.062.42'
$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("D19").Value = '2.083.38'
$ws.Range("E19").Value = '  +0.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '221.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9992'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.164'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.56%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9971'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1403'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.527'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.89'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.496'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.185'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.91%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05626'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.42%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.146'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.203'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7462'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.845'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.36%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.141'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.20%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.663'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.39%  '

$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.248.85'
$ws.Range("E38").Value = '  +0.75%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.767'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01779'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.53%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.581'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.32%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8947'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9985'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.06%  '

$ws.Range("D45").Value = '1.984.62'
$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.10%  '

$ws.Range("E47").Value = '  -1.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5067'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.97%  '

$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4074'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.19%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.018'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.45%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05836'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.74%  '
